# ---------------------------------------------------------------------------
# Applies the "Adjusted data input and crude hyper parameter tuning for svm"
# commit to the workbook:
#   * rename "test" sheet to "test_results"
#   * convert D2:D17 "concat" formulas into a proper shared formula
#   * add a new "train_results" sheet with NN/RandomForest custom-vs-tfidf
#     training accuracies, plus a bar chart on that sheet
#   * point the four existing charts on the "Charts" sheet at the renamed
#     "test_results" sheet, retitle them, and clamp chart1's value axis
#   * restore the various cell selections that were left behind by the
#     original author
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename "test" -> "test_results"
# ---------------------------------------------------------------------------
$wsTest = $wb.Worksheets.Item("test")
$wsTest.Name = "test_results"

# Rebuild D2:D17 as a single shared formula (matches t="shared" ref/si markup)
$wsTest.Range("D2:D17").Formula = "=A2&""_""&B2&""_""&C2"

# Restore the selection that was left on the sheet
$wsTest.Range("D21").Select()

# ---------------------------------------------------------------------------
# 2. Update the four charts living on the "Charts" sheet
# ---------------------------------------------------------------------------
$wsCharts = $wb.Worksheets.Item("Charts")

$chart1 = $wsCharts.ChartObjects().Item(1).Chart
$chart1.ChartTitle.Text = "Accuracy"
$chart1.SeriesCollection(1).Formula = "=SERIES(,test_results!`$D`$2:`$D`$17,test_results!`$E`$2:`$E`$17,1)"
$chart1.Axes(2).MinimumScale = 0.5

$chart2 = $wsCharts.ChartObjects().Item(2).Chart
$chart2.ChartTitle.Text = "f_positive"
$chart2.SeriesCollection(1).Formula = "=SERIES(,test_results!`$D`$2:`$D`$17,test_results!`$F`$2:`$F`$17,1)"

$chart3 = $wsCharts.ChartObjects().Item(3).Chart
$chart3.ChartTitle.Text = "f_neutral"
$chart3.SeriesCollection(1).Formula = "=SERIES(,test_results!`$D`$2:`$D`$17,test_results!`$I`$2:`$I`$17,1)"

$chart4 = $wsCharts.ChartObjects().Item(4).Chart
$chart4.ChartTitle.Text = "f_negative"
$chart4.SeriesCollection(1).Formula = "=SERIES(,test_results!`$D`$2:`$D`$17,test_results!`$L`$2:`$L`$17,1)"

# ---------------------------------------------------------------------------
# 3. Add the new "train_results" sheet (after "Charts") with training scores
# ---------------------------------------------------------------------------
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$wsTrain = $wb.Worksheets.Add($null, $lastSheet)
$wsTrain.Name = "train_results"

$wsTrain.Range("A2").Value = "NN_TFIDF"
$wsTrain.Range("A1").Value = "NN_Custom"
$wsTrain.Range("A4").Value = "RandomForest_TFIDF"
$wsTrain.Range("A3").Value = "RandomForest_Custom"
$wsTrain.Range("B1").Value = 0.99
$wsTrain.Range("B2").Value = 0.901
$wsTrain.Range("B3").Value = 0.992
$wsTrain.Range("B4").Value = 0.917

# Bar chart comparing custom vs tf-idf features on the training set
$trainChartObj = $wsTrain.ChartObjects().Add(300, 300, 400, 300)
$trainChartObj.Chart.ChartType = 51
$trainChartObj.Chart.SetSourceData($wsTrain.Range("A1:B4"))
$trainChartObj.Chart.HasTitle = $true
$trainChartObj.Chart.ChartTitle.Text = "Custom V TF-IDF in Train"

$wsTrain.Range("A5").Select()

# ---------------------------------------------------------------------------
# 4. Restore remaining selections / active sheet
# ---------------------------------------------------------------------------
$wsCharts.Range("L13").Select()
$wsCharts.Activate()
